# all_avg_num_ridesV2.xlsx -- "Add files via upload"
#
# The re-uploaded workbook normalises the member_casual labels in column B
# for the older rows (2-71) from the capitalised "Casual"/"Member" strings
# to the same lowercase "casual"/"member" strings already used by every
# later row, and leaves the UI selection parked near the bottom of the
# data (around I220) instead of I7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalise every "Casual"/"Member" value in column B (member_casual) to
# the lowercase spelling used elsewhere in the sheet. Walking the whole
# used range (instead of a hard-coded row count) makes this resilient to
# exactly which rows still have the old capitalised text.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $cur = $cell.Value2
    if ($cur -eq "Casual") {
        $cell.Value = "casual"
    } elseif ($cur -eq "Member") {
        $cell.Value = "member"
    }
}

# Match the saved view state: active cell/selection moved from I7 down to
# I220 (near the end of the data), with the window scrolled so row 219 is
# at the top.
$ws.Range("I220").Select() | Out-Null
try { $excel.ActiveWindow.ScrollRow = 219 } catch {}
